$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.920.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.52%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.967.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.18%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.79%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.50%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.687"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.80%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.795"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.38%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.13%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.29%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000331"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.52%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.25%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.599.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.17%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.965.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.22%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.11%  "

# Row 17
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.70%  "

# Row 18
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.43%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.838.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.51%  "

# Row 20
$ws.Range("E20").Value = "  -0.91%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "452.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.90%  "

# Row 22
$ws.Range("E22").Value = "  +5.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.78%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.47%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.28%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.93%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.66%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.37%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.19%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.43%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.08%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "14.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.54%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.68%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.129"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.60%  "

# Row 35
$ws.Range("E35").Value = "  +14.83%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "69.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.49%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "637.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.16%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.431"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.48%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.62%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.147"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.34%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.01%  "

# Row 42
$ws.Range("E42").Value = "  +0.04%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +45.03%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0484"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.69%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.58"
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.149"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.81%  "

# Row 47
$ws.Range("E47").Value = "  -10.76%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.00%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000290"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.87%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.76%  "

# Row 51
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.825.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.01%  "
